# Update "想去人数" (attendance count) figures on the 展览 and 全部类型 sheets
# F2: 1548 -> 1549
# F3: 79   -> 80

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1549
    $ws.Range("F3").Value = 80
}
